$wb = $excel.ActiveWorkbook

# --- Rename the "Include from RoleCode" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from RoleCode")
$wsInclude.Name = "Include #0"

# --- Update metadata values on the "Metadata" sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update Version value
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" after "Contact" (row 10) and before "Description" (row 11)
$wsMeta.Rows.Item(11).Insert()

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# Copy the formatting of the row above (Contact, row 10) onto the new row so it matches the rest of the table
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
